$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (instead of Excel auto-converting numeric-looking
# strings like prices/percentages into numbers) so the values round-trip
# verbatim, matching the original inline-string cell content.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "274.91"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.61%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.24"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.88%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.883"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.02%"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.24%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.947"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.76%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.210"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-1.52%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8758"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.37%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1512"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "3.83%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.05121"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.13%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07536"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.37%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02955"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.79%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08982"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.65%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001574"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.29%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006375"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.15%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006195"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2.48%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.466"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.33%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.312"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.25%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.55%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.95%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.76%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.919"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.87%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04410"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.07%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001179"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.30%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.003855"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-12.42%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001201"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.00%"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "14.64%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04126"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.41%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006779"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.71%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.83%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002091"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.47%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01146"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-2.39%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005186"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-2.30%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-28.10%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.02001"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.08%"
